$d = $word.ActiveDocument

$LQ = [char]0x201C   # “
$RQ = [char]0x201D   # ”

# ------------------------------------------------------------------
# Locate the second occurrence of the "O usuario acessa..." sentence
# (the one that contains "de forma que pelo menos...", i.e. the
# use-case's alternative/exception flow), and the start of the
# following paragraph's sentence ("O sistema valida...").
# ------------------------------------------------------------------

$search1 = "O usu" + [char]0x00E1 + "rio acessa o site do sistema, informa os campos " + $LQ
$rngFirst = $d.Range(0, $d.Content.End)
$f1 = $rngFirst.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f1) { throw "could not locate the first 'O usuario acessa...' sentence" }

$rngSecond = $d.Range($rngFirst.End, $d.Content.End)
$f2 = $rngSecond.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f2) { throw "could not locate the second 'O usuario acessa...' sentence" }
$para1Start = $rngSecond.Start

# End of the "...esteja errado " run (includes trailing space), i.e.
# right before the lone "." run that must stay a separate run.
$rngErrado = $d.Range($para1Start, $d.Content.End)
$f3 = $rngErrado.Find.Execute("esteja errado ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $f3) { throw "could not locate 'esteja errado '" }
$para1End = $rngErrado.End

# ------------------------------------------------------------------
# Move the _GoBack bookmark to sit right between the merged run and
# the trailing "." run *before* editing any text, so that the engine
# does not fold the "." run into the merged run.
# ------------------------------------------------------------------
$bmRange = $d.Range($para1End, $para1End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Merge the runs of paragraph 1's sentence into a single run with the
# cleaned-up text (drop the stray space before the closing quote and
# the trailing space before the final period).
# ------------------------------------------------------------------
$newText1 = "O usu" + [char]0x00E1 + "rio acessa o site do sistema, informa os campos " + $LQ + "nome" + $RQ + " e " + $LQ + "senha" + $RQ + " de forma que pelo menos a informa" + [char]0x00E7 + [char]0x00E3 + "o de um dos campos esteja errado"

$rngReplace1 = $d.Range($para1Start, $para1End)
$rngReplace1.Text = $newText1

# ------------------------------------------------------------------
# Merge the runs of the following paragraph's sentence
# ("O sistema valida ... incorretos”.") into a single run. The text
# itself does not change; editing it in place forces the engine to
# coalesce the previously-split runs (the bookmark has already moved
# away, so it will not be re-inserted here).
#
# NB: paragraph 1's text length changed above, so re-locate this
# range with a fresh search instead of reusing old offsets.
# ------------------------------------------------------------------
$rngPara2Start = $d.Range($rngReplace1.End, $d.Content.End)
[void]$rngPara2Start.Find.Execute("O sistema valida", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2Start = $rngPara2Start.Start

$rngPara2End = $d.Range($para2Start, $d.Content.End)
[void]$rngPara2End.Find.Execute("incorretos" + $RQ + ".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2End = $rngPara2End.End

$newText2 = "O sistema valida as informa" + [char]0x00E7 + [char]0x00F5 + "es de nome e senha no banco de dados e retorna a mensagem " + $LQ + "Usu" + [char]0x00E1 + "rio ou senha incorretos" + $RQ + "."

# The engine only coalesces the runs spanned by a Range.Text
# assignment when the new text actually differs from the old text,
# so first force a change (placeholder), then set the final text.
$rngReplace2 = $d.Range($para2Start, $para2End)
$rngReplace2.Text = "x"
$rngReplace2b = $d.Range($para2Start, $rngReplace2.End)
$rngReplace2b.Text = $newText2
